# Auto-generated: apply scheduled-runner market data refresh to Hyperion Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2561
$ws.Range("J18").Value = 2492.5715
$ws.Range("L18").Value = 2492.5715
$ws.Range("N18").Value = -3060.5715
$ws.Range("H32").Value = 2256.5557
$ws.Range("I32").Value = 1199.75
$ws.Range("J32").Value = 3102
$ws.Range("K32").Value = 1199.75
$ws.Range("L32").Value = 3102
$ws.Range("M32").Value = -873.75
$ws.Range("N32").Value = -3754
$ws.Range("H40").Value = 5073.3335
$ws.Range("I40").Value = 5179.7144
$ws.Range("J40").Value = 4701
$ws.Range("K40").Value = 5179.7144
$ws.Range("L40").Value = 4701
$ws.Range("M40").Value = -5004.7144
$ws.Range("N40").Value = -5051
$ws.Range("H43").Value = 1489.6428
$ws.Range("I43").Value = 1613.3334
$ws.Range("K43").Value = 1613.3334
$ws.Range("M43").Value = -1544.3334
$ws.Range("H51").Value = 4514.0713
$ws.Range("I51").Value = 3727.818
$ws.Range("J51").Value = 7397
$ws.Range("K51").Value = 3727.818
$ws.Range("L51").Value = 7397
$ws.Range("M51").Value = -3243.818
$ws.Range("N51").Value = -8365
$ws.Range("H127").Value = 2373.6667
$ws.Range("I127").Value = 2051.625
$ws.Range("K127").Value = 6154.875
$ws.Range("M127").Value = -1194.875
$ws.Range("H134").Value = 190758
$ws.Range("J134").Value = 190758
$ws.Range("L134").Value = 190758
$ws.Range("N134").Value = -200898

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 62500300
$ws.Range("I5").Value = 83333700
$ws.Range("J5").Value = 99
$ws.Range("K5").Value = 83333700
$ws.Range("L5").Value = 99
$ws.Range("M5").Value = -83333588
$ws.Range("N5").Value = -323
$ws.Range("H32").Value = 6259.0894
$ws.Range("I32").Value = 4634.675
$ws.Range("K32").Value = 4634.675
$ws.Range("M32").Value = -4347.675
$ws.Range("H45").Value = 7995321
$ws.Range("I45").Value = 10277092
$ws.Range("K45").Value = 10277092
$ws.Range("M45").Value = -10276715
$ws.Range("H88").Value = 1700.375
$ws.Range("I88").Value = 274.5
$ws.Range("J88").Value = 2175.6667
$ws.Range("K88").Value = 274.5
$ws.Range("L88").Value = 2175.6667
$ws.Range("M88").Value = 131.5
$ws.Range("N88").Value = -2987.6667
$ws.Range("H91").Value = 1700.375
$ws.Range("I91").Value = 274.5
$ws.Range("J91").Value = 2175.6667
$ws.Range("K91").Value = 274.5
$ws.Range("L91").Value = 2175.6667
$ws.Range("M91").Value = 1129.5
$ws.Range("N91").Value = -4983.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 62500300
$ws.Range("I4").Value = 83333700
$ws.Range("J4").Value = 99
$ws.Range("K4").Value = 83333700
$ws.Range("L4").Value = 99
$ws.Range("M4").Value = -83333585
$ws.Range("N4").Value = -329
$ws.Range("H11").Value = 200
$ws.Range("I11").Value = 200
$ws.Range("K11").Value = 200
$ws.Range("M11").Value = -60
$ws.Range("H22").Value = 2645769.5
$ws.Range("J22").Value = 202
$ws.Range("L22").Value = 202
$ws.Range("N22").Value = -548
$ws.Range("H86").Value = 3130741.8
$ws.Range("I86").Value = 3852482
$ws.Range("J86").Value = 3200.3333
$ws.Range("K86").Value = 3852482
$ws.Range("L86").Value = 3200.3333
$ws.Range("M86").Value = -3851359
$ws.Range("N86").Value = -5446.3333
$ws.Range("H89").Value = 3130741.8
$ws.Range("I89").Value = 3852482
$ws.Range("J89").Value = 3200.3333
$ws.Range("K89").Value = 19262410
$ws.Range("L89").Value = 16001.6665
$ws.Range("M89").Value = -19256794
$ws.Range("N89").Value = -27233.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H7").Value = 319.83334
$ws.Range("I7").Value = 99.28570999999999
$ws.Range("K7").Value = 99.28570999999999
$ws.Range("M7").Value = 13.71429000000001
$ws.Range("H16").Value = 2596.4167
$ws.Range("J16").Value = 3348.5
$ws.Range("L16").Value = 3348.5
$ws.Range("N16").Value = -3922.5
$ws.Range("H23").Value = 9133.333000000001
$ws.Range("I23").Value = 8000
$ws.Range("K23").Value = 8000
$ws.Range("M23").Value = -7760
$ws.Range("H27").Value = 9133.333000000001
$ws.Range("I27").Value = 8000
$ws.Range("K27").Value = 8000
$ws.Range("M27").Value = -7808
$ws.Range("H55").Value = 14512.167
$ws.Range("I55").Value = 10691
$ws.Range("J55").Value = 18333.334
$ws.Range("K55").Value = 10691
$ws.Range("L55").Value = 18333.334
$ws.Range("M55").Value = -10376
$ws.Range("N55").Value = -18963.334
$ws.Range("H105").Value = 4372.5
$ws.Range("I105").Value = 4000
$ws.Range("J105").Value = 4745
$ws.Range("K105").Value = 4000
$ws.Range("L105").Value = 4745
$ws.Range("M105").Value = -2253
$ws.Range("N105").Value = -8239
$ws.Range("H113").Value = 2596.4167
$ws.Range("J113").Value = 3348.5
$ws.Range("L113").Value = 3348.5
$ws.Range("N113").Value = -7688.5
$ws.Range("H132").Value = 91950.28999999999
$ws.Range("I132").Value = 64128.688
$ws.Range("K132").Value = 192386.064
$ws.Range("M132").Value = -189856.064

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 4779.8
$ws.Range("I11").Value = 20000
$ws.Range("J11").Value = 974.75
$ws.Range("K11").Value = 60000
$ws.Range("L11").Value = 2924.25
$ws.Range("M11").Value = -59860
$ws.Range("N11").Value = -3204.25
$ws.Range("H26").Value = 291.5
$ws.Range("J26").Value = 92
$ws.Range("L26").Value = 276
$ws.Range("N26").Value = -852
$ws.Range("H61").Value = 137
$ws.Range("I61").Value = 137
$ws.Range("K61").Value = 411
$ws.Range("M61").Value = -196
$ws.Range("H109").Value = 1592
$ws.Range("I109").Value = 1024
$ws.Range("J109").Value = 5000
$ws.Range("K109").Value = 3072
$ws.Range("L109").Value = 15000
$ws.Range("M109").Value = -2032
$ws.Range("N109").Value = -17080

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6152.875
$ws.Range("I7").Value = 2806.25
$ws.Range("J7").Value = 9499.5
$ws.Range("K7").Value = 2806.25
$ws.Range("L7").Value = 9499.5
$ws.Range("M7").Value = -2694.25
$ws.Range("N7").Value = -9723.5
$ws.Range("H22").Value = 124385.375
$ws.Range("I22").Value = 165581.5
$ws.Range("J22").Value = 797
$ws.Range("K22").Value = 165581.5
$ws.Range("L22").Value = 797
$ws.Range("M22").Value = -165286.5
$ws.Range("N22").Value = -1387
$ws.Range("H27").Value = 124385.375
$ws.Range("I27").Value = 165581.5
$ws.Range("J27").Value = 797
$ws.Range("K27").Value = 165581.5
$ws.Range("L27").Value = 797
$ws.Range("M27").Value = -165474.5
$ws.Range("N27").Value = -1011
$ws.Range("H126").Value = 6152.875
$ws.Range("I126").Value = 2806.25
$ws.Range("J126").Value = 9499.5
$ws.Range("K126").Value = 8418.75
$ws.Range("L126").Value = 28498.5
$ws.Range("M126").Value = -5948.75
$ws.Range("N126").Value = -33438.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9059.84
$ws.Range("J62").Value = 9647.695
$ws.Range("L62").Value = 9647.695
$ws.Range("N62").Value = -10895.695
$ws.Range("H65").Value = 9059.84
$ws.Range("J65").Value = 9647.695
$ws.Range("L65").Value = 48238.475
$ws.Range("N65").Value = -54478.475
